$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts all existing
# data (previously columns A:BY) one column to the right (now B:BZ) and
# keeps all existing content/formatting intact.
$ws.Columns.Item(1).Insert()

# --- New column A: section labels --------------------------------------
$ws.Range("A9").Value  = "Piyush Comments"
$ws.Range("A10").Value = "Ekta Commnets"

# --- Row 10: additional "combine levels" comments ----------------------
$ws.Range("G10").Value = "Combine IR2 and IR3 only"
$ws.Range("Q10").Value = "Combine Hip and Mansard"
$ws.Range("R10").Value = "Combine ClyTile, Membran, Metal, Roll"
$ws.Range("S10").Value = "Cluster"
$ws.Range("T10").Value = "Cluster"

# --- Rows 13 & 14: extra notes area -------------------------------------
$ws.Range("AJ14").Value = "combine Floor and OthW"
$ws.Range("AK14").Value = "combine Fa and Po"
$ws.Range("AM14").Value = "combine Mix and FuseP"

# --- Styling ------------------------------------------------------------
# Row 9 (A9:T9) and row 10 (A10:T10) switch to the plain-black "Calibri"
# font used by the new annotation area, except I9 which keeps its
# original red "note" font.
$ws.Range("A9:H9").Font.Color = 0
$ws.Range("J9:T9").Font.Color = 0
$ws.Range("A10:T10").Font.Color = 0

# Rows 13 and 14 (V:CA) are touched so they become part of the used
# range, using the same plain black font (no fill).
$ws.Range("V13:CA14").Font.Color = 0
